$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column before D, shifting the existing D:K (now E:L)
#        one column to the right. Seed the freshly-inserted column D with
#        column E's values/format as a starting point (most rows in this
#        report repeat the same figure/blank/"NA" across every period), then
#        overwrite the cells that actually carry a new reporting period's
#        figures below.
$ws.Columns("D").Insert()
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Populate the brand-new column D with the latest reporting
#        period's figures (a new fiscal year was added in front of
#        the existing history).
$newColumnD = @{
    7   = 43465
    8   = 632000
    15  = -14400
    17  = 101400
    18  = 530600
    20  = -322000
    21  = 246900
    23  = 208700
    24  = 17900
    26  = 190800
    27  = 190800
    29  = 0
    32  = 322000
    33  = 190800
    35  = 190800
    38  = 43465
    41  = 284000
    42  = 542500
    48  = 485900
    49  = 1113300
    52  = 90300
    54  = 19728400
    58  = 632600
    59  = 181000
    61  = 247700
    66  = 17038900
    72  = 527700
    76  = 2689600
    80  = 43465
    81  = 190800
    83  = 38200
    89  = 234400
    91  = -33400
    94  = -271400
    96  = -82200
    100 = 63700
    102 = 26700
}

foreach ($row in $newColumnD.Keys) {
    $ws.Cells.Item($row, 4).Value = $newColumnD[$row]
}
